$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '67.305.76'
$ws.Range('E2').Value = '  -1.18%  '

$ws.Range('D3').Value = '3.517.09'
$ws.Range('E3').Value = '  +0.07%  '

$ws.Range('E4').Value = '  +0.06%  '

$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '611.38'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +0.37%  '

$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '148.46'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  -1.79%  '

$ws.Range('D7').Value = '3.514.47'
$ws.Range('E7').Value = '  +0.04%  '

$ws.Range('E8').Value = '  +0.13%  '

$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.481'
$ws.Range('D9').Style = 'Normal'

$ws.Range('E10').Value = '  -1.52%  '

$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '8.06'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  +6.71%  '

$ws.Range('E12').Value = '  -1.92%  '

$ws.Range('E13').Value = '  +0.37%  '

$ws.Range('D14').Value = '4.113.49'
$ws.Range('E14').Value = '  +0.25%  '

$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '31.58'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  -1.68%  '

$ws.Range('D16').Value = '3.515.14'
$ws.Range('E16').Value = '  -0.19%  '

$ws.Range('D17').Value = '67.339.15'
$ws.Range('E17').Value = '  -1.14%  '

$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '0.117'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  -0.05%  '

$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '10.92'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  +9.29%  '

$ws.Range('E20').Value = '  -2.50%  '

$ws.Range('E21').Value = '  -0.15%  '

$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '437.10'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  -3.20%  '

$ws.Range('E23').Value = '  -3.03%  '

$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '80.11'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  +1.28%  '

$ws.Range('D25').Value = '3.658.66'
$ws.Range('E25').Value = '  +0.32%  '

$ws.Range('E26').Value = '  +0.07%  '

$ws.Range('E27').Value = '  -4.32%  '

$ws.Range('E28').Value = '  -1.65%  '

$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '8.28'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  -5.21%  '

$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '2.53'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  +0.27%  '

$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '1.59'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  -5.23%  '

$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '1.00'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  +0.01%  '

$ws.Range('E33').Value = '  -2.44%  '

$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '25.61'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  -0.38%  '

$ws.Range('E35').Value = '  -3.74%  '

$ws.Range('E36').Value = '  -2.14%  '

$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '8.05'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  -0.02%  '

$ws.Range('E39').Value = '  +0.03%  '

$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '176.49'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  -0.36%  '

$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.0901'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  -0.67%  '

$ws.Range('E42').Value = '  -0.75%  '

$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '2.06'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  -10.38%  '

$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.897'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  -0.44%  '

$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '46.38'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  -1.29%  '

$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '28.19'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  -8.55%  '

$ws.Range('E47').Value = '  -5.29%  '

$ws.Range('E48').Value = '  -2.17%  '

$ws.Range('E49').Value = '  -2.35%  '

$ws.Range('E50').Value = '  -1.08%  '

$ws.Range('E51').Value = '  -2.42%  '
